$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 50, shifting the old rows 50-56 down to 51-57.
# This splits the combined "Le rapport...\nURL" citation (previously on two rows,
# 50 and 51) so that the row layout matches the new source citation block:
#   49 Source:
#   50 (blank)
#   51 Le rapport de l'enquete sur les entreprises a Madagascar, INSTAT, p. 9
#   52 (blank)
#   53 http://www.instat.mg/pdf/enq_entreprises_2005.pdf
#   56 NISR
#   57 <new MTI citation>
$ws.Rows.Item(50).Insert()

# After the insert, the old A51 (which held the hyperlinked URL text, style
# "HyperLink") is now at A52. The new layout wants A52 to hold the same blank
# "source"-styled cell pattern used at A50/A52 elsewhere, so restore plain
# (non-hyperlink) formatting there.
$ws.Range("A52").Font.Underline = $false
$ws.Range("A52").Font.Color = $ws.Range("A49").Font.Color

# Drop the hyperlink entirely (target text relocates to a NEW, separate cell
# A53 as plain text, with no hyperlink attached).
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $ws.Hyperlinks.Item($i).Delete()
}

# Update the final citation text (was the NISR Establishment Census citation,
# now the MTI SME Development Policy citation). This cell is now at A57.
$ws.Range("A57").Value = 'Ministry of Trade and Industry (MTI), "Small and Medium Enterprises (SMEs) Development Policy", 2010, p. 7. Available at http://www.rdb.rw/uploads/tx_sbdownloader/SME_Devt_policy_V180610.pdf'
